$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.061289673249482
$ws.Range("D2").Value = 1.06334133136097
$ws.Range("E2").Value = 1.065701894906011
$ws.Range("F2").Value = 1.075239413303145
$ws.Range("I2").Value = 1.049826041592673
$ws.Range("J2").Value = 1.066266120853881
$ws.Range("K2").Value = 1.066060232648585
$ws.Range("L2").Value = 1.06841442467493
$ws.Range("M2").Value = 1.077926507773063
$ws.Range("N2").Value = 1.067780341047927
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.062431550793067
$ws.Range("D3").Value = 1.064237074861781
$ws.Range("E3").Value = 1.066764942224118
$ws.Range("F3").Value = 1.076322527149728
$ws.Range("I3").Value = 1.05015300399774
$ws.Range("J3").Value = 1.067061065156079
$ws.Range("K3").Value = 1.066770765359522
$ws.Range("L3").Value = 1.069292303421557
$ws.Range("M3").Value = 1.078826247068643
$ws.Range("N3").Value = 1.068576414262216
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.063170372880071
$ws.Range("D4").Value = 1.064816671534027
$ws.Range("E4").Value = 1.067453613139572
$ws.Range("F4").Value = 1.077023947183265
$ws.Range("I4").Value = 1.050363373319048
$ws.Range("J4").Value = 1.067574818728567
$ws.Range("K4").Value = 1.067229881563803
$ws.Range("L4").Value = 1.069860553497489
$ws.Range("M4").Value = 1.079408409468476
$ws.Range("N4").Value = 1.069090897423705
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.063480962763616
$ws.Range("D5").Value = 1.065060331695015
$ws.Range("E5").Value = 1.067743323792951
$ws.Range("F5").Value = 1.077318961365095
$ws.Range("I5").Value = 1.050451526027969
$ws.Range("J5").Value = 1.067790650517856
$ws.Range("K5").Value = 1.067422739471729
$ws.Range("L5").Value = 1.07009949459184
$ws.Range("M5").Value = 1.079653143637339
$ws.Range("N5").Value = 1.069307035718891
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.063533111555148
$ws.Range("D6").Value = 1.065101243172754
$ws.Range("E6").Value = 1.067791978825037
$ws.Range("F6").Value = 1.077368503590872
$ws.Range("I6").Value = 1.050466310453901
$ws.Range("J6").Value = 1.067826880824139
$ws.Range("K6").Value = 1.067455112103551
$ws.Range("L6").Value = 1.070139616711498
$ws.Range("M6").Value = 1.079694235178108
$ws.Range("N6").Value = 1.069343317476365
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.063174523039175
$ws.Range("D7").Value = 1.064819927342057
$ws.Range("E7").Value = 1.067457483507674
$ws.Range("F7").Value = 1.077027888638051
$ws.Range("I7").Value = 1.050364552345737
$ws.Range("J7").Value = 1.067577703273398
$ws.Range("K7").Value = 1.067232459147731
$ws.Range("L7").Value = 1.069863746048522
$ws.Range("M7").Value = 1.079411679645467
$ws.Range("N7").Value = 1.069093786064921
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.061675586756974
$ws.Range("D8").Value = 1.063644053953921
$ws.Range("E8").Value = 1.066060989324535
$ws.Range("F8").Value = 1.075605338086518
$ws.Range("I8").Value = 1.049936787997669
$ws.Range("J8").Value = 1.066534906305698
$ws.Range("K8").Value = 1.066300494319357
$ws.Range("L8").Value = 1.068711065683941
$ws.Range("M8").Value = 1.078230584939332
$ws.Range("N8").Value = 1.068049508205915
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.059033857048312
$ws.Range("D9").Value = 1.0615719459908
$ws.Range("E9").Value = 1.063606381569974
$ws.Range("F9").Value = 1.073103008319199
$ws.Range("I9").Value = 1.049173840296874
$ws.Range("J9").Value = 1.06469253874912
$ws.Range("K9").Value = 1.064653300860348
$ws.Range("L9").Value = 1.066681457731416
$ws.Range("M9").Value = 1.076149119264627
$ws.Range("N9").Value = 1.066204524276085
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.057272372391009
$ws.Range("D10").Value = 1.06019049452921
$ws.Range("E10").Value = 1.061974142297669
$ws.Range("F10").Value = 1.071437731476608
$ws.Range("I10").Value = 1.048659036016301
$ws.Range("J10").Value = 1.063461029522197
$ws.Range("K10").Value = 1.063551832864676
$ws.Range("L10").Value = 1.065329434895976
$ws.Range("M10").Value = 1.074761316520546
$ws.Range("N10").Value = 1.064971266164804
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.056509538559242
$ws.Range("D11").Value = 1.059592298951793
$ws.Range("E11").Value = 1.061268352217289
$ws.Range("F11").Value = 1.07071734397297
$ws.Range("I11").Value = 1.048434654363792
$ws.Range("J11").Value = 1.062926993110511
$ws.Range("K11").Value = 1.063074090258595
$ws.Range("L11").Value = 1.064744241831402
$ws.Range("M11").Value = 1.074160342471171
$ws.Range("N11").Value = 1.064436471360158
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.056226171965229
$ws.Range("D12").Value = 1.059370099595584
$ws.Range("E12").Value = 1.061006337485092
$ws.Range("F12").Value = 1.070449863092731
$ws.Range("I12").Value = 1.048351088265265
$ws.Range("J12").Value = 1.062728509764535
$ws.Range("K12").Value = 1.062896514898537
$ws.Range("L12").Value = 1.064526911201669
$ws.Range("M12").Value = 1.073937106682967
$ws.Range("N12").Value = 1.064237706145065
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.056286955811617
$ws.Range("D13").Value = 1.059417762245222
$ws.Range("E13").Value = 1.061062533879075
$ws.Range("F13").Value = 1.070507233988286
$ws.Range("I13").Value = 1.048369023470405
$ws.Range("J13").Value = 1.06277109048453
$ws.Range("K13").Value = 1.062934610873691
$ws.Range("L13").Value = 1.064573527727558
$ws.Range("M13").Value = 1.073984991851208
$ws.Range("N13").Value = 1.064280347334566
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.056486115724541
$ws.Range("D14").Value = 1.059573931943588
$ws.Range("E14").Value = 1.061246690997263
$ws.Range("F14").Value = 1.070695231818233
$ws.Range("I14").Value = 1.048427751268718
$ws.Range("J14").Value = 1.062910588834563
$ws.Range("K14").Value = 1.063059414283748
$ws.Range("L14").Value = 1.064726276477798
$ws.Range("M14").Value = 1.074141889881339
$ws.Range("N14").Value = 1.064420043788257
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.056608822544433
$ws.Range("D15").Value = 1.05967015286135
$ws.Range("E15").Value = 1.061360175767296
$ws.Range("F15").Value = 1.07081107714048
$ws.Range("I15").Value = 1.048463906144091
$ws.Range("J15").Value = 1.062996522660874
$ws.Range("K15").Value = 1.063136293808189
$ws.Range("L15").Value = 1.06482039480612
$ws.Range("M15").Value = 1.07423855896363
$ws.Range("N15").Value = 1.064506099650458
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.057322996962918
$ws.Range("D16").Value = 1.060230194405178
$ws.Range("E16").Value = 1.062021003901775
$ws.Range("F16").Value = 1.071485555699623
$ws.Range("I16").Value = 1.048673896540393
$ws.Range("J16").Value = 1.063496455179393
$ws.Range("K16").Value = 1.063583522190379
$ws.Range("L16").Value = 1.065368277309498
$ws.Range("M16").Value = 1.074801200186791
$ws.Range("N16").Value = 1.065006742130497
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.057770952248643
$ws.Range("D17").Value = 1.060581488727417
$ws.Range("E17").Value = 1.062435786279216
$ws.Range("F17").Value = 1.071908822358677
$ws.Range("I17").Value = 1.048805224792672
$ws.Range("J17").Value = 1.063809839149259
$ws.Range("K17").Value = 1.063863842398134
$ws.Range("L17").Value = 1.065712014117361
$ws.Range("M17").Value = 1.075154117544813
$ws.Range("N17").Value = 1.065320571141545
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.058032227483618
$ws.Range("D18").Value = 1.060786391039833
$ws.Range("E18").Value = 1.062677816416738
$ws.Range("F18").Value = 1.072155773076793
$ws.Range("I18").Value = 1.048881684713542
$ws.Range("J18").Value = 1.063992555032126
$ws.Range("K18").Value = 1.064027271296415
$ws.Range("L18").Value = 1.065912533319393
$ws.Range("M18").Value = 1.07535996382214
$ws.Range("N18").Value = 1.065503546501923
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.058121314011735
$ws.Range("D19").Value = 1.060856257146635
$ws.Range("E19").Value = 1.062760358488398
$ws.Range("F19").Value = 1.072239988168014
$ws.Range("I19").Value = 1.048907731554755
$ws.Range("J19").Value = 1.064054843605345
$ws.Range("K19").Value = 1.064082983254032
$ws.Range("L19").Value = 1.065980909149628
$ws.Range("M19").Value = 1.075430151338402
$ws.Range("N19").Value = 1.065565923532061
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.057722891877569
$ws.Range("D20").Value = 1.060543798350603
$ws.Range("E20").Value = 1.062391274272201
$ws.Range("F20").Value = 1.071863402982088
$ws.Range("I20").Value = 1.048791149164665
$ws.Range("J20").Value = 1.063776223847312
$ws.Range("K20").Value = 1.063833774673668
$ws.Range("L20").Value = 1.065675131998312
$ws.Range("M20").Value = 1.075116253314343
$ws.Range("N20").Value = 1.065286908102013
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.056427468549525
$ws.Range("D21").Value = 1.059527943940933
$ws.Range("E21").Value = 1.061192457255667
$ws.Range("F21").Value = 1.070639868299026
$ws.Range("I21").Value = 1.048410463491673
$ws.Range("J21").Value = 1.062869513321358
$ws.Range("K21").Value = 1.063022666118817
$ws.Range("L21").Value = 1.06468129478216
$ws.Range("M21").Value = 1.07409568752517
$ws.Range("N21").Value = 1.064378909943112
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.055612890484438
$ws.Range("D22").Value = 1.058889219176921
$ws.Range("E22").Value = 1.060439564585921
$ws.Range("F22").Value = 1.069871179916834
$ws.Range("I22").Value = 1.048169833934209
$ws.Range("J22").Value = 1.062298742875905
$ws.Range("K22").Value = 1.0625119927229
$ws.Range("L22").Value = 1.064056638824085
$ws.Range("M22").Value = 1.073453975104108
$ws.Range("N22").Value = 1.063807328938158
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.056044723045554
$ws.Range("D23").Value = 1.059227820911061
$ws.Range("E23").Value = 1.060838606604066
$ws.Range("F23").Value = 1.070278619738622
$ws.Range("I23").Value = 1.048297517365223
$ws.Range("J23").Value = 1.062601384250427
$ws.Range("K23").Value = 1.062782776557281
$ws.Range("L23").Value = 1.064387761023814
$ws.Range("M23").Value = 1.073794163104317
$ws.Range("N23").Value = 1.064110400098145
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.057744608319511
$ws.Range("D24").Value = 1.06056082901468
$ws.Range("E24").Value = 1.062411387039736
$ws.Range("F24").Value = 1.07188392583784
$ws.Range("I24").Value = 1.048797509772575
$ws.Range("J24").Value = 1.063791413388933
$ws.Range("K24").Value = 1.063847361222262
$ws.Range("L24").Value = 1.065691797368074
$ws.Range("M24").Value = 1.075133362543799
$ws.Range("N24").Value = 1.065302119214526
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.059716861967966
$ws.Range("D25").Value = 1.062107643706108
$ws.Range("E25").Value = 1.064240222375449
$ws.Range("F25").Value = 1.073749401386022
$ws.Range("I25").Value = 1.049372168103513
$ws.Range("J25").Value = 1.065169408838429
$ws.Range("K25").Value = 1.065079727694248
$ws.Range("L25").Value = 1.067205975036686
$ws.Range("M25").Value = 1.076687255534348
$ws.Range("N25").Value = 1.066682071575618
Write-Host "Updated 264 cells across rows 2-25 (columns B-N, excluding G and H)"
